{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that separates it from the \"LOQ4073...\" line)\n// that the site-generator used to append after the \"Requisitos\" section.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the \"Ver no Jupiter...\" paragraph; the footer block is that paragraph,\n// the \"\u00a9 2020 ...\" paragraph right after it, and the blank paragraph right\n// before it (the spacer that used to separate it from \"LOQ4073...\").\nlet verIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetTexts[0]) {\n    verIndex = i;\n    break;\n  }\n}\n\nif (verIndex !== -1) {\n  const toDelete = [];\n  // Blank spacer paragraph right before \"Ver no Jupiter...\"\n  if (verIndex - 1 >= 0 && paragraphs.items[verIndex - 1].text === \"\") {\n    toDelete.push(paragraphs.items[verIndex - 1]);\n  }\n  // \"Ver no Jupiter...\" itself\n  toDelete.push(paragraphs.items[verIndex]);\n  // \"\u00a9 2020 ...\" right after it\n  if (\n    verIndex + 1 < paragraphs.items.length &&\n    paragraphs.items[verIndex + 1].text === targetTexts[1]\n  ) {\n    toDelete.push(paragraphs.items[verIndex + 1]);\n  }\n\n  // Delete from the end backwards so earlier indices stay valid.\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph that separates it from the \"LOQ4073...\" line)\n# that the site-generator used to append after the \"Requisitos\" section.\n\n$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"Contact: luizeleno@usp.br. Powered by Jekyll and Github pages.\"\n\n$count = $d.Paragraphs.Count\n$verIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*$verText*\") {\n        $verIndex = $i\n        break\n    }\n}\n\nif ($verIndex -ge 1) {\n    $toDelete = New-Object System.Collections.ArrayList\n\n    # The paragraph right after \"Ver no Jupiter...\" -- the \"(c) 2020 ...\" line\n    if (($verIndex + 1) -le $count) {\n        $nextText = $d.Paragraphs.Item($verIndex + 1).Range.Text\n        if ($nextText -like \"*$copyrightText*\") {\n            [void]$toDelete.Add($verIndex + 1)\n        }\n    }\n\n    # \"Ver no Jupiter...\" itself\n    [void]$toDelete.Add($verIndex)\n\n    # The blank spacer paragraph right before \"Ver no Jupiter...\"\n    if (($verIndex - 1) -ge 1) {\n        $prevText = $d.Paragraphs.Item($verIndex - 1).Range.Text\n        if ($prevText.Trim().Length -eq 0) {\n            [void]$toDelete.Add($verIndex - 1)\n        }\n    }\n\n    # Delete from the highest index down so earlier indices stay valid.\n    $sorted = $toDelete | Sort-Object -Descending\n    foreach ($idx in $sorted) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
